# Atualização automática de QUARAI.xlsx
#
# Changes:
#  1. Rename sheet "Paineis DARQ" -> "PAINEIS DARQ"
#  2. Rename sheet "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#  3. Delete sheet "Desarquivamentos Pendentes" (no longer needed)

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true
